# BOQ_Upload_Format.xlsx -- Works and Finance changes
# Adds a "Scope/Milestone" column at the front of the table, renames the
# remaining headers, and appends four more sample/milestone rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new column before A to hold the Scope/Milestone values.
#    (Everything that used to be in A:E shifts right to B:F.)
# ---------------------------------------------------------------------
$ws.Columns("A:A").Insert()

# ---------------------------------------------------------------------
# 2. Header row - copy the existing header look onto the new A1 cell,
#    then fix up the text for every header cell.
# ---------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("A1").Value = "Scope/Milestone"
$ws.Range("B1").Value = "Item Description"
$ws.Range("C1").Value = "Ref DSR/NS"
$ws.Range("D1").Value = "Unit"
$ws.Range("E1").Value = "Rate"
$ws.Range("F1").Value = "Quantity"

# ---------------------------------------------------------------------
# 3. Fill in the Scope/Milestone values for the existing data rows
#    (2-4), and give that column a boxed border.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "MileStone1"
$ws.Range("A3").Value = "MileStone2"
$ws.Range("A4").Value = "MileStone3"

# ---------------------------------------------------------------------
# 4. Add the new rows 5-8 - copy the formatting of row 4 down first so
#    the new rows inherit the same look, then overwrite the values.
# ---------------------------------------------------------------------
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F8").PasteSpecial(-4122)

$ws.Range("A5").Value = "MileStone3"
$ws.Range("B5").Value = "Test455"
$ws.Range("C5").Value = "Test3 Desc gou"
$ws.Range("D5").Value = "cm"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 5

$ws.Range("A6").Value = "MileStone4"
$ws.Range("B6").Value = "Test456"
$ws.Range("C6").Value = "Test3 Desc gou"
$ws.Range("D6").Value = "cm"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 5

$ws.Range("A7").Value = "MileStone5"
$ws.Range("B7").Value = "Test457"
$ws.Range("C7").Value = "Test3 Desc gou"
$ws.Range("D7").Value = "cm"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 5

$ws.Range("A8").Value = "MileStone5"
$ws.Range("B8").Value = "Test458"
$ws.Range("C8").Value = "Test3 Desc gou"
$ws.Range("D8").Value = "cm"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 5

# ---------------------------------------------------------------------
# 5. Box-border the whole Scope/Milestone column (header + data).
# ---------------------------------------------------------------------
$ws.Range("A1:A8").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 6. Auto-fit the new column, restore the explicit width afterwards so
#    it matches the rest of the sheet's "customWidth" columns.
# ---------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 14.375

# ---------------------------------------------------------------------
# 7. Row-height bookkeeping at the bottom of the sheet: the two blank
#    formatted rows that used to trail the sheet (999/1000) go away,
#    and rows 19/20 (just past the new table) pick up the same blank
#    15.75pt row height used by the rest of the sheet.
# ---------------------------------------------------------------------
$ws.Rows("999:1000").Delete()
$ws.Rows("19:20").RowHeight = 15.75

# ---------------------------------------------------------------------
# 8. Selection, matching the saved cursor position in the workbook.
# ---------------------------------------------------------------------
$ws.Range("C14").Select()
